# Updated cryptos list on Sun Aug 18 10:11:18 UTC 2024 with GitHub Actions
#
# Refreshes the "Coin" (B), "Link" (C), "Price" (D) and "Volume(1h)" (E)
# columns of the crypto table on the active sheet, row by row, to match
# the latest coinranking.com snapshot. A handful of rows also re-sort
# (ShibaInu <-> WrappedEther, Maker <-> InjectiveProtocol <-> RenderToken)
# so their Coin/Link text moves along with the row's new Price/Volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel coercing
# numeric-looking strings (e.g. "1.00", "6.60") into real numbers and
# silently dropping the formatting. We force it in as text (quote-prefix
# trick, same as typing '1.00 into the cell) and then reset the cell
# style back to Normal so no stray quote-prefix formatting is left behind.
function Set-TextValue($cell, [string]$value) {
    $looksNumeric = $value -match '^[+-]?[0-9]+(\.[0-9]+)?$'
    if ($looksNumeric) {
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

# row -> @{ col = newValue; ... }  (only cells that actually changed)
$updates = [ordered]@{
    2  = @{ D = "59.910.36";  E = "  +1.44%  " }
    3  = @{ D = "2.639.58";   E = "  +1.91%  " }
    4  = @{ D = "1.00";       E = "  -0.29%  " }
    5  = @{ D = "538.13";     E = "  +2.14%  " }
    6  = @{ D = "143.55";     E = "  +2.97%  " }
    7  = @{ D = "0.998";      E = "  -0.03%  " }
    8  = @{                   E = "  +0.77%  " }
    9  = @{ D = "6.60";       E = "  +1.32%  " }
    10 = @{                   E = "  +1.78%  " }
    11 = @{ D = "0.337";      E = "  +1.41%  " }
    12 = @{                   E = "  -1.54%  " }
    13 = @{ D = "3.101.73";   E = "  +1.69%  " }
    14 = @{ D = "59.817.75";  E = "  +1.39%  " }
    15 = @{ D = "20.91";      E = "  +1.86%  " }
    16 = @{ B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "2.654.57";   E = "  +2.47%  " }
    17 = @{ B = "ShibaInu";      C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib";                 D = "0.0000134"; E = "  +1.13%  " }
    18 = @{ D = "343.25";     E = "  -0.01%  " }
    19 = @{                   E = "  +1.51%  " }
    20 = @{ D = "10.20";      E = "  +1.01%  " }
    21 = @{ D = "6.41";       E = "  -0.31%  " }
    23 = @{ D = "67.44";      E = "  +1.62%  " }
    24 = @{ D = "0.412";      E = "  +1.51%  " }
    25 = @{ D = "0.167";      E = "  -1.13%  " }
    26 = @{ D = "0.998";      E = "  -0.05%  " }
    27 = @{ D = "7.27";       E = "  +2.83%  " }
    28 = @{ D = "0.0₃0753";   E = "  +4.15%  " }
    29 = @{                   E = "  -0.02%  " }
    30 = @{                   E = "  +3.80%  " }
    31 = @{ D = "5.86";       E = "  -1.15%  " }
    32 = @{ D = "18.91";      E = "  +0.98%  " }
    33 = @{ D = "151.08";     E = "  +1.25%  " }
    34 = @{ D = "4.01";       E = "  +1.01%  " }
    35 = @{                   E = "  +1.00%  " }
    36 = @{                   E = "  -1.07%  " }
    37 = @{ D = "0.839";      E = "  +3.51%  " }
    38 = @{ D = "0.838";      E = "  +1.36%  " }
    39 = @{ D = "3.56";       E = "  +0.94%  " }
    40 = @{ D = "281.99";     E = "  +3.99%  " }
    41 = @{ D = "0.998";      E = "  +0.06%  " }
    42 = @{                   E = "  +0.70%  " }
    44 = @{ D = "0.0951";     E = "  -0.23%  " }
    45 = @{                   E = "  +2.73%  " }
    46 = @{ B = "Maker";              C = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr";         D = "1.972.11"; E = "  +0.40%  " }
    47 = @{ B = "InjectiveProtocol";  C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D = "18.66";    E = "  +2.15%  " }
    48 = @{                   E = "  +1.18%  " }
    49 = @{ B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "4.55"; E = "  -1.71%  " }
    50 = @{ D = "112.55";     E = "  -1.68%  " }
    51 = @{ D = "4.74";       E = "  +0.50%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        Set-TextValue $cell $cols[$col]
    }
}
